$wb = $excel.ActiveWorkbook

# Update Estadisticos Matutinos data (Aprobados, Reprobados, Por_Apro, Por_Repro,
# Promedio, Blancos, Por_Blan columns E:K) for the affected "Docente" rows.
# The "1er Parcial" and "3er Parcial" sheets hold identical statistic tables,
# so the same updates are applied to both; "2o Parcial" is left untouched.
$targetSheetNames = @("1er Parcial", "3er Parcial")

foreach ($sheetName in $targetSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(5, 5).Value = 17
    $ws.Cells.Item(5, 6).Value = 5
    $ws.Cells.Item(5, 7).Value = 77.27
    $ws.Cells.Item(5, 8).Value = 22.73
    $ws.Cells.Item(5, 9).Value = 6.9
    $ws.Cells.Item(5, 10).Value = 5
    $ws.Cells.Item(5, 11).Value = 22.73
    $ws.Cells.Item(6, 5).Value = 27
    $ws.Cells.Item(6, 6).Value = 4
    $ws.Cells.Item(6, 7).Value = 87.09999999999999
    $ws.Cells.Item(6, 8).Value = 12.9
    $ws.Cells.Item(6, 9).Value = 8.5
    $ws.Cells.Item(6, 10).Value = 4
    $ws.Cells.Item(6, 11).Value = 12.9
    $ws.Cells.Item(7, 5).Value = 22
    $ws.Cells.Item(7, 6).Value = 16
    $ws.Cells.Item(7, 7).Value = 57.89
    $ws.Cells.Item(7, 8).Value = 42.11
    $ws.Cells.Item(7, 9).Value = 8.4
    $ws.Cells.Item(7, 10).Value = 14
    $ws.Cells.Item(7, 11).Value = 36.84
    $ws.Cells.Item(8, 5).Value = 24
    $ws.Cells.Item(8, 6).Value = 11
    $ws.Cells.Item(8, 7).Value = 68.56999999999999
    $ws.Cells.Item(8, 8).Value = 31.43
    $ws.Cells.Item(8, 9).Value = 7.9
    $ws.Cells.Item(8, 10).Value = 10
    $ws.Cells.Item(8, 11).Value = 28.57
    $ws.Cells.Item(9, 5).Value = 21
    $ws.Cells.Item(9, 6).Value = 12
    $ws.Cells.Item(9, 7).Value = 63.64
    $ws.Cells.Item(9, 8).Value = 36.36
    $ws.Cells.Item(9, 9).Value = 7.8
    $ws.Cells.Item(9, 10).Value = 9
    $ws.Cells.Item(9, 11).Value = 27.27
    $ws.Cells.Item(13, 5).Value = 22
    $ws.Cells.Item(13, 6).Value = 2
    $ws.Cells.Item(13, 7).Value = 91.67
    $ws.Cells.Item(13, 8).Value = 8.33
    $ws.Cells.Item(13, 9).Value = 6.7
    $ws.Cells.Item(13, 10).Value = 2
    $ws.Cells.Item(13, 11).Value = 8.33
    $ws.Cells.Item(14, 5).Value = 24
    $ws.Cells.Item(14, 6).Value = 7
    $ws.Cells.Item(14, 7).Value = 77.42
    $ws.Cells.Item(14, 8).Value = 22.58
    $ws.Cells.Item(14, 9).Value = 8
    $ws.Cells.Item(14, 10).Value = 7
    $ws.Cells.Item(14, 11).Value = 22.58
    $ws.Cells.Item(15, 5).Value = 11
    $ws.Cells.Item(15, 6).Value = 10
    $ws.Cells.Item(15, 7).Value = 52.38
    $ws.Cells.Item(15, 8).Value = 47.62
    $ws.Cells.Item(15, 9).Value = 8
    $ws.Cells.Item(15, 10).Value = 10
    $ws.Cells.Item(15, 11).Value = 47.62
    $ws.Cells.Item(16, 5).Value = 28
    $ws.Cells.Item(16, 6).Value = 7
    $ws.Cells.Item(16, 7).Value = 80
    $ws.Cells.Item(16, 8).Value = 20
    $ws.Cells.Item(16, 9).Value = 8.6
    $ws.Cells.Item(16, 10).Value = 7
    $ws.Cells.Item(16, 11).Value = 20
    $ws.Cells.Item(17, 5).Value = 13
    $ws.Cells.Item(17, 6).Value = 8
    $ws.Cells.Item(17, 7).Value = 61.9
    $ws.Cells.Item(17, 8).Value = 38.1
    $ws.Cells.Item(17, 9).Value = 7.8
    $ws.Cells.Item(17, 10).Value = 8
    $ws.Cells.Item(17, 11).Value = 38.1
    $ws.Cells.Item(25, 5).Value = 17
    $ws.Cells.Item(25, 6).Value = 22
    $ws.Cells.Item(25, 7).Value = 43.59
    $ws.Cells.Item(25, 8).Value = 56.41
    $ws.Cells.Item(25, 9).Value = 8.300000000000001
    $ws.Cells.Item(25, 10).Value = 22
    $ws.Cells.Item(25, 11).Value = 56.41
    $ws.Cells.Item(26, 5).Value = 15
    $ws.Cells.Item(26, 6).Value = 20
    $ws.Cells.Item(26, 7).Value = 42.86
    $ws.Cells.Item(26, 8).Value = 57.14
    $ws.Cells.Item(26, 9).Value = 8.1
    $ws.Cells.Item(26, 10).Value = 20
    $ws.Cells.Item(26, 11).Value = 57.14
    $ws.Cells.Item(27, 5).Value = 13
    $ws.Cells.Item(27, 6).Value = 19
    $ws.Cells.Item(27, 7).Value = 40.63
    $ws.Cells.Item(27, 8).Value = 59.38
    $ws.Cells.Item(27, 9).Value = 7.2
    $ws.Cells.Item(27, 10).Value = 19
    $ws.Cells.Item(27, 11).Value = 59.38
    $ws.Cells.Item(28, 5).Value = 29
    $ws.Cells.Item(28, 6).Value = 14
    $ws.Cells.Item(28, 7).Value = 67.44
    $ws.Cells.Item(28, 8).Value = 32.56
    $ws.Cells.Item(28, 9).Value = 7.4
    $ws.Cells.Item(28, 10).Value = 14
    $ws.Cells.Item(28, 11).Value = 32.56
    $ws.Cells.Item(29, 5).Value = 20
    $ws.Cells.Item(29, 6).Value = 9
    $ws.Cells.Item(29, 7).Value = 68.97
    $ws.Cells.Item(29, 8).Value = 31.03
    $ws.Cells.Item(29, 9).Value = 7.7
    $ws.Cells.Item(29, 10).Value = 9
    $ws.Cells.Item(29, 11).Value = 31.03
    $ws.Cells.Item(30, 5).Value = 24
    $ws.Cells.Item(30, 6).Value = 10
    $ws.Cells.Item(30, 7).Value = 70.59
    $ws.Cells.Item(30, 8).Value = 29.41
    $ws.Cells.Item(30, 9).Value = 7.6
    $ws.Cells.Item(30, 10).Value = 0
    $ws.Cells.Item(30, 11).Value = 0
    $ws.Cells.Item(31, 5).Value = 21
    $ws.Cells.Item(31, 6).Value = 10
    $ws.Cells.Item(31, 7).Value = 67.73999999999999
    $ws.Cells.Item(31, 8).Value = 32.26
    $ws.Cells.Item(31, 9).Value = 7.3
    $ws.Cells.Item(31, 10).Value = 0
    $ws.Cells.Item(31, 11).Value = 0
    $ws.Cells.Item(34, 5).Value = 21
    $ws.Cells.Item(34, 6).Value = 12
    $ws.Cells.Item(34, 7).Value = 63.64
    $ws.Cells.Item(34, 8).Value = 36.36
    $ws.Cells.Item(34, 9).Value = 6.6
    $ws.Cells.Item(34, 10).Value = 0
    $ws.Cells.Item(34, 11).Value = 0
}
